$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new data rows right before the totals row (old row 10) ---
# This pushes the old "totals" row (10) down to 12, and the footer row
# (old 11) down to 13 - mirroring the structural insert in the diff.
$ws.Rows.Item(10).EntireRow.Insert()
$ws.Rows.Item(10).EntireRow.Insert()

# --- Clone formatting + merges of the row-9 item row onto the two new rows ---
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial()
$ws.Range("A9:Q9").Copy()
$ws.Range("A11:Q11").PasteSpecial()

# --- Row 10: new item "ORGASOL LIGHT CREAM" ---
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "ORGASOL LIGHT CREAM"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "0"
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "130.00"
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "130.0000"
$ws.Range("Q10").NumberFormat = "@"
$ws.Range("Q10").Value = "1:0"

# --- Row 11: new item "PRISBRINA  CAPS" ---
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "PRISBRINA  CAPS"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "-1:-1"
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = "0"
$ws.Range("N11").NumberFormat = "@"
$ws.Range("N11").Value = "150.00"
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "150.0000"
$ws.Range("Q11").NumberFormat = "@"
$ws.Range("Q11").Value = "1:0"

# --- Update the running total (old P10 -> now shifted to P12) ---
$ws.Range("P12").Value = 656.5
